$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (never let Excel auto-coerce
# numeric-looking strings like "6.42" or "1.00" into real numbers), and
# reset the cell style back to the workbook default afterwards so we don't
# leave a stray NumberFormat/quote-prefix style on the cell.
function Set-CellText($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-CellText 2 4 '53.986.66'
Set-CellText 2 5 '  +1.52%  '
# Row 3 - Ethereum
Set-CellText 3 4 '2.245.85'
Set-CellText 3 5 '  +3.05%  '
# Row 4 - TetherUSD
Set-CellText 4 4 '1.00'
Set-CellText 4 5 '  -0.02%  '
# Row 5 - BNB
Set-CellText 5 4 '494.33'
Set-CellText 5 5 '  +3.12%  '
# Row 6 - Solana
Set-CellText 6 4 '127.53'
Set-CellText 6 5 '  +3.55%  '
# Row 7 - USDC
Set-CellText 7 4 '0.996'
Set-CellText 7 5 '  -0.21%  '
# Row 8 - XRP
Set-CellText 8 5 '  +2.45%  '
# Row 9 - Dogecoin
Set-CellText 9 4 '0.0950'
Set-CellText 9 5 '  +4.88%  '
# Row 10 - TRON
Set-CellText 10 5 '  +2.88%  '
# Row 11 - Cardano
Set-CellText 11 5 '  +4.76%  '
# Row 12 - Toncoin
Set-CellText 12 5 '  +1.85%  '
# Row 13 - WrappedliquidstakedEther2.0
Set-CellText 13 4 '2.649.78'
Set-CellText 13 5 '  +3.22%  '
# Row 14 - Avalanche
Set-CellText 14 4 '21.73'
Set-CellText 14 5 '  +4.25%  '
# Row 15 - WrappedBTC
Set-CellText 15 4 '53.921.64'
Set-CellText 15 5 '  +1.55%  '
# Row 16 - ShibaInu
Set-CellText 16 4 '0.0000129'
Set-CellText 16 5 '  +2.13%  '
# Row 17 - WrappedEther
Set-CellText 17 4 '2.259.34'
Set-CellText 17 5 '  +3.99%  '
# Row 18 - Chainlink
Set-CellText 18 4 '10.01'
Set-CellText 18 5 '  +6.05%  '
# Row 19 - Polkadot
Set-CellText 19 5 '  +4.79%  '
# Row 20 - was BitcoinCash, now Uniswap (rows 20/21 swapped)
Set-CellText 20 2 'Uniswap'
Set-CellText 20 3 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-CellText 20 4 '6.42'
Set-CellText 20 5 '  +6.96%  '
# Row 21 - was Uniswap, now BitcoinCash
Set-CellText 21 2 'BitcoinCash'
Set-CellText 21 3 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-CellText 21 4 '299.39'
Set-CellText 21 5 '  +2.57%  '
# Row 22 - Dai
Set-CellText 22 4 '0.997'
Set-CellText 22 5 '  -0.11%  '
# Row 23 - LEO
Set-CellText 23 4 '5.37'
Set-CellText 23 5 '  -3.39%  '
# Row 24 - Litecoin
Set-CellText 24 4 '61.89'
Set-CellText 24 5 '  -0.20%  '
# Row 25 - Binance-PegBSC-USD
Set-CellText 25 4 '1.02'
Set-CellText 25 5 '  +2.21%  '
# Row 26 - Polygon
Set-CellText 26 5 '  +2.69%  '
# Row 27 - WrappedeETH
Set-CellText 27 4 '2.360.03'
Set-CellText 27 5 '  +3.61%  '
# Row 28 - Kaspa
Set-CellText 28 5 '  +3.82%  '
# Row 29 - InternetComputer(DFINITY)
Set-CellText 29 4 '7.04'
Set-CellText 29 5 '  +1.55%  '
# Row 30 - Monero
Set-CellText 30 4 '166.53'
Set-CellText 30 5 '  +0.64%  '
# Row 31 - PancakeSwap
Set-CellText 31 5 '  +2.81%  '
# Row 32 - PEPE
Set-CellText 32 4 '0.0₃0683'
Set-CellText 32 5 '  +5.16%  '
# Row 33 - Aptos
Set-CellText 33 5 '  +4.48%  '
# Row 34 - USDe
Set-CellText 34 5 '  +0.04%  '
# Row 35 - FirstDigitalUSD
Set-CellText 35 5 '  -0.15%  '
# Row 36 - Fetch.AI
Set-CellText 36 5 '  +2.91%  '
# Row 37 - EthereumClassic
Set-CellText 37 4 '17.63'
Set-CellText 37 5 '  +3.06%  '
# Row 38 - SuiNetwork
Set-CellText 38 4 '0.903'
Set-CellText 38 5 '  +11.99%  '
# Row 39 - ImmutableX
Set-CellText 39 5 '  +4.83%  '
# Row 40 - NEARProtocol
Set-CellText 40 5 '  +4.94%  '
# Row 41 - OKB
Set-CellText 41 4 '35.68'
Set-CellText 41 5 '  +0.63%  '
# Row 42 - Stacks
Set-CellText 42 5 '  +4.29%  '
# Row 43 - PolygonEcosystemToken
Set-CellText 43 4 '0.371'
Set-CellText 43 5 '  +2.78%  '
# Row 44 - Filecoin
Set-CellText 44 4 '3.35'
Set-CellText 44 5 '  +4.24%  '
# Row 45 - was Aave, now RenderToken (rows 45/46 swapped)
Set-CellText 45 2 'RenderToken'
Set-CellText 45 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-CellText 45 4 '4.93'
Set-CellText 45 5 '  +5.55%  '
# Row 46 - was RenderToken, now Aave
Set-CellText 46 2 'Aave'
Set-CellText 46 3 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-CellText 46 4 '125.45'
Set-CellText 46 5 '  +2.87%  '
# Row 47 - Stellar
Set-CellText 47 4 '0.0886'
Set-CellText 47 5 '  +1.93%  '
# Row 48 - Mantle
Set-CellText 48 4 '0.540'
Set-CellText 48 5 '  +2.84%  '
# Row 49 - Bittensor
Set-CellText 49 4 '236.22'
Set-CellText 49 5 '  +4.43%  '
# Row 50 - Hedera
Set-CellText 50 4 '0.0483'
Set-CellText 50 5 '  +4.20%  '
# Row 51 - VeChain
Set-CellText 51 4 '0.0202'
Set-CellText 51 5 '  +2.30%  '
